$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 74; this shifts existing rows 74-93 down to 75-94
$ws.Rows.Item(74).Insert()

# Populate the new row 74 with the latest weekly record
$ws.Cells.Item(74, 1).Value = 11
$ws.Cells.Item(74, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(74, 3).Value = "Bíobío"
$ws.Cells.Item(74, 4).Value = 44551
$ws.Cells.Item(74, 5).Value = 8
$ws.Cells.Item(74, 6).Value = 100112043
$ws.Cells.Item(74, 7).Value = "Pepino ensalada"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 100
$ws.Cells.Item(74, 11).Value = 7000
$ws.Cells.Item(74, 12).Value = 8000
$ws.Cells.Item(74, 13).Value = 7500
$ws.Cells.Item(74, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(74, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(74, 16).Value = 125
$ws.Cells.Item(74, 17).Value = 60
$ws.Cells.Item(74, 18).Value = "Hortaliza"
